# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Old = 4513; New = 4530 },
    @{ Row = 3;  Old = 2490; New = 2494 },
    @{ Row = 8;  Old = 223;  New = 225 },
    @{ Row = 10; Old = 168;  New = 171 },
    @{ Row = 11; Old = 170;  New = 171 },
    @{ Row = 12; Old = 1683; New = 1693 },
    @{ Row = 13; Old = 304;  New = 308 },
    @{ Row = 14; Old = 3674; New = 3698 },
    @{ Row = 15; Old = 17;   New = 21 },
    @{ Row = 16; Old = 245;  New = 246 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.New
}

$updates2 = @(
    @{ Row = 2;  Old = 4513; New = 4530 },
    @{ Row = 3;  Old = 2490; New = 2494 },
    @{ Row = 10; Old = 223;  New = 225 },
    @{ Row = 12; Old = 168;  New = 171 },
    @{ Row = 13; Old = 170;  New = 171 },
    @{ Row = 16; Old = 1683; New = 1693 },
    @{ Row = 17; Old = 304;  New = 308 },
    @{ Row = 18; Old = 3675; New = 3698 },
    @{ Row = 19; Old = 17;   New = 21 },
    @{ Row = 20; Old = 245;  New = 246 }
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates2) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.New
}
